$wb = $excel.ActiveWorkbook

# The same updates apply to both the "展览" and "全部类型" sheets:
#   F3: 5601 -> 5602
#   F5: 680  -> 681
#   F7: 27   -> 28
#   F15: 247 -> 250
#   F19: 4577 -> 4581

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 5602
    $ws.Range("F5").Value = 681
    $ws.Range("F7").Value = 28
    $ws.Range("F15").Value = 250
    $ws.Range("F19").Value = 4581
}
